$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 69
$ws.Cells.Item(69, 1).Value = 111785228
$ws.Cells.Item(69, 2).Value = 78578
$ws.Cells.Item(69, 4).Value = "NT"
$ws.Cells.Item(69, 5).Value = 6458
$ws.Cells.Item(69, 6).Value = "Lunglav"
$ws.Cells.Item(69, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(69, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(69, 17).Value = 577256
$ws.Cells.Item(69, 18).Value = 6944531
$ws.Cells.Item(69, 26).ClearContents()
$ws.Cells.Item(69, 28).ClearContents()

# Row 70
$ws.Cells.Item(70, 1).Value = 111785230
$ws.Cells.Item(70, 2).Value = 78578
$ws.Cells.Item(70, 4).Value = "NT"
$ws.Cells.Item(70, 5).Value = 6458
$ws.Cells.Item(70, 6).Value = "Lunglav"
$ws.Cells.Item(70, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(70, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(70, 17).Value = 577262
$ws.Cells.Item(70, 18).Value = 6944620
$ws.Cells.Item(70, 26).ClearContents()
$ws.Cells.Item(70, 28).ClearContents()

# Row 71
$ws.Cells.Item(71, 1).Value = 111785200
$ws.Cells.Item(71, 2).Value = 78512
$ws.Cells.Item(71, 4).Value = "LC"
$ws.Cells.Item(71, 5).Value = 6456
$ws.Cells.Item(71, 6).Value = "Skinnlav"
$ws.Cells.Item(71, 7).Value = "Leptogium saturninum"
$ws.Cells.Item(71, 8).Value = "(Dicks.) Nyl."
$ws.Cells.Item(71, 17).Value = 577256
$ws.Cells.Item(71, 18).Value = 6944531
$ws.Cells.Item(71, 26).ClearContents()
$ws.Cells.Item(71, 28).ClearContents()

# Row 72
$ws.Cells.Item(72, 1).Value = 111785244
$ws.Cells.Item(72, 2).Value = 96348
$ws.Cells.Item(72, 4).Value = "VU"
$ws.Cells.Item(72, 5).Value = 220787
$ws.Cells.Item(72, 6).Value = "Knärot"
$ws.Cells.Item(72, 7).Value = "Goodyera repens"
$ws.Cells.Item(72, 8).Value = "(L.) R. Br."
$ws.Cells.Item(72, 17).Value = 577364
$ws.Cells.Item(72, 18).Value = 6944622
$ws.Cells.Item(72, 26).ClearContents()
$ws.Cells.Item(72, 28).ClearContents()

# Row 73
$ws.Cells.Item(73, 1).Value = 111785235
$ws.Cells.Item(73, 2).Value = 77267
$ws.Cells.Item(73, 4).Value = "NT"
$ws.Cells.Item(73, 5).Value = 6446
$ws.Cells.Item(73, 6).Value = "Kolflarnlav"
$ws.Cells.Item(73, 7).Value = "Carbonicola anthracophila"
$ws.Cells.Item(73, 8).Value = "(Nyl.) Bendiksby & Timdal"
$ws.Cells.Item(73, 17).Value = 577227
$ws.Cells.Item(73, 18).Value = 6944649
$ws.Cells.Item(73, 26).ClearContents()
$ws.Cells.Item(73, 28).ClearContents()

# Row 74
$ws.Cells.Item(74, 1).Value = 111785251
$ws.Cells.Item(74, 2).Value = 93161
$ws.Cells.Item(74, 4).Value = "VU"
$ws.Cells.Item(74, 5).Value = 1079
$ws.Cells.Item(74, 6).Value = "Aspfjädermossa"
$ws.Cells.Item(74, 7).Value = "Neckera pennata"
$ws.Cells.Item(74, 8).Value = "Hedw."
$ws.Cells.Item(74, 17).Value = 577283
$ws.Cells.Item(74, 18).Value = 6944534
$ws.Cells.Item(74, 26).ClearContents()
$ws.Cells.Item(74, 28).ClearContents()

# Row 75
$ws.Cells.Item(75, 1).Value = 111785201
$ws.Cells.Item(75, 2).Value = 78512
$ws.Cells.Item(75, 4).Value = "LC"
$ws.Cells.Item(75, 5).Value = 6456
$ws.Cells.Item(75, 6).Value = "Skinnlav"
$ws.Cells.Item(75, 7).Value = "Leptogium saturninum"
$ws.Cells.Item(75, 8).Value = "(Dicks.) Nyl."
$ws.Cells.Item(75, 17).Value = 577248
$ws.Cells.Item(75, 18).Value = 6944531
$ws.Cells.Item(75, 26).ClearContents()
$ws.Cells.Item(75, 28).ClearContents()

# Row 76
$ws.Cells.Item(76, 1).Value = 111785191
$ws.Cells.Item(76, 2).Value = 89405
$ws.Cells.Item(76, 4).Value = "NT"
$ws.Cells.Item(76, 5).Value = 1202
$ws.Cells.Item(76, 6).Value = "Ullticka"
$ws.Cells.Item(76, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(76, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(76, 17).Value = 577236
$ws.Cells.Item(76, 18).Value = 6944656
$ws.Cells.Item(76, 26).ClearContents()
$ws.Cells.Item(76, 28).ClearContents()

# Row 77
$ws.Cells.Item(77, 1).Value = 111785199
$ws.Cells.Item(77, 2).Value = 89416
$ws.Cells.Item(77, 4).Value = "LC"
$ws.Cells.Item(77, 5).Value = 1205
$ws.Cells.Item(77, 6).Value = "Stor aspticka"
$ws.Cells.Item(77, 7).Value = "Phellinus populicola"
$ws.Cells.Item(77, 8).Value = "Niemelä"
$ws.Cells.Item(77, 17).Value = 577256
$ws.Cells.Item(77, 18).Value = 6944531
$ws.Cells.Item(77, 26).ClearContents()
$ws.Cells.Item(77, 28).ClearContents()

# Row 78
$ws.Cells.Item(78, 1).Value = 111785192
$ws.Cells.Item(78, 2).Value = 89405
$ws.Cells.Item(78, 4).Value = "NT"
$ws.Cells.Item(78, 5).Value = 1202
$ws.Cells.Item(78, 6).Value = "Ullticka"
$ws.Cells.Item(78, 7).Value = "Phellinidium ferrugineofuscum"
$ws.Cells.Item(78, 8).Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Cells.Item(78, 17).Value = 577282
$ws.Cells.Item(78, 18).Value = 6944714
$ws.Cells.Item(78, 26).ClearContents()
$ws.Cells.Item(78, 28).ClearContents()

# Row 79
$ws.Cells.Item(79, 1).Value = 111785202
$ws.Cells.Item(79, 2).Value = 78512
$ws.Cells.Item(79, 4).Value = "LC"
$ws.Cells.Item(79, 5).Value = 6456
$ws.Cells.Item(79, 6).Value = "Skinnlav"
$ws.Cells.Item(79, 7).Value = "Leptogium saturninum"
$ws.Cells.Item(79, 8).Value = "(Dicks.) Nyl."
$ws.Cells.Item(79, 17).Value = 577215
$ws.Cells.Item(79, 18).Value = 6944631
$ws.Cells.Item(79, 26).ClearContents()
$ws.Cells.Item(79, 28).ClearContents()

# Row 80
$ws.Cells.Item(80, 1).Value = 111785190
$ws.Cells.Item(80, 2).Value = 94134
$ws.Cells.Item(80, 4).Value = "NT"
$ws.Cells.Item(80, 5).Value = 53
$ws.Cells.Item(80, 6).Value = "Vedtrappmossa"
$ws.Cells.Item(80, 7).Value = "Crossocalyx hellerianus"
$ws.Cells.Item(80, 8).Value = "(Nees ex Lindenb.) Meyl."
$ws.Cells.Item(80, 17).Value = 577243
$ws.Cells.Item(80, 18).Value = 6944541
$ws.Cells.Item(80, 26).ClearContents()
$ws.Cells.Item(80, 28).ClearContents()

# Row 81
$ws.Cells.Item(81, 1).Value = 111785229
$ws.Cells.Item(81, 2).Value = 78578
$ws.Cells.Item(81, 4).Value = "NT"
$ws.Cells.Item(81, 5).Value = 6458
$ws.Cells.Item(81, 6).Value = "Lunglav"
$ws.Cells.Item(81, 7).Value = "Lobaria pulmonaria"
$ws.Cells.Item(81, 8).Value = "(L.) Hoffm."
$ws.Cells.Item(81, 17).Value = 577208
$ws.Cells.Item(81, 18).Value = 6944522
$ws.Cells.Item(81, 26).ClearContents()
$ws.Cells.Item(81, 28).ClearContents()

# Row 82
$ws.Cells.Item(82, 1).Value = 111785206
$ws.Cells.Item(82, 2).Value = 77268
$ws.Cells.Item(82, 4).Value = "NT"
$ws.Cells.Item(82, 5).Value = 228912
$ws.Cells.Item(82, 6).Value = "Mörk kolflarnlav"
$ws.Cells.Item(82, 7).Value = "Carbonicola myrmecina"
$ws.Cells.Item(82, 8).Value = "(Ach.) Bendiksby & Timdal"
$ws.Cells.Item(82, 17).Value = 577236
$ws.Cells.Item(82, 18).Value = 6944656
$ws.Cells.Item(82, 26).ClearContents()
$ws.Cells.Item(82, 28).ClearContents()
